$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("valid_login").Name = "test_valid_login"
$wb.Worksheets.Item("invalid_login").Name = "test_invalid_login"

# Update the active selection on the invalid_login sheet (now renamed)
$ws = $wb.Worksheets.Item("test_invalid_login")
$ws.Activate()
$ws.Range("D27").Select()
